# "Mi avance de automation"
# DriverConfig: point the Browser row's link at the new endpoint and drop the
# now-unused second (Firefox / old QA url) row; move selection to B3.

$wb = $excel.ActiveWorkbook
$driverConfig = $wb.Worksheets.Item("DriverConfig")

# --- keep B2's existing look (Hyperlink style, centered/wrapped, etc.) -----
# The engine's Hyperlinks collection only supports a full-sheet Delete, so we
# stash B2's current formatting in a scratch cell, rebuild the hyperlink, and
# paste the formatting back afterwards.
$driverConfig.Range("Z1").Value = "fmt-backup"
$driverConfig.Range("B2").Copy() | Out-Null
$driverConfig.Range("Z1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Remove every hyperlink on the sheet (B2's + the stale B3 one) ...
$driverConfig.Hyperlinks.Delete()

# ... then rebuild only the one that should survive, now aimed at the new URL.
$driverConfig.Hyperlinks.Add($driverConfig.Range("B2"), "https://hl.com") | Out-Null
$driverConfig.Range("B2").Value = "https://hl.com"

$driverConfig.Range("Z1").Copy() | Out-Null
$driverConfig.Range("B2").PasteSpecial(-4122) | Out-Null   # restore original format
$driverConfig.Range("Z1").Clear()

$excel.CutCopyMode = 0

# --- drop the Firefox row's contents (keep the row's own cell styles) -----
$driverConfig.Range("A3").ClearContents()
$driverConfig.Range("B3").ClearContents()

# --- match the saved selection (B3 active on the DriverConfig tab) --------
$driverConfig.Activate()
$driverConfig.Range("B3").Select()
